$wb = $excel.ActiveWorkbook

# --- "Spell Modifiers" sheet: add a new "Bouncing Projectile" modifier row ---
$wsModifiers = $wb.Worksheets.Item("Spell Modifiers")
$wsModifiers.Range("A6").Value = "Bouncing Projectile"
$wsModifiers.Range("B6").Value = "Ricochets off walls"

# --- "OnHit Effect" sheet: expand the Explosion effect's description ---
$wsOnHit = $wb.Worksheets.Item("OnHit Effect")
$wsOnHit.Range("B11").Value = "Explodes. Yep. Applies other effects in spell."

# --- Update selections / active sheet to match the final view state ---
$wsModifiers.Activate()
$wsModifiers.Range("B7").Select()

$wsOnHit.Activate()
$wsOnHit.Range("B12").Select()
